$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing the existing "HOLY FAMILY MARONITE CHURCH"
# record down to row 17 (carrying its per-column formatting: text style for
# A/B/C/E, date style for D).
$ws.Rows.Item(16).Insert()

# The inserted row reverts to the sheet's default row height; restore the
# 13.05pt custom height used by every other data row.
$ws.Rows.Item(16).RowHeight = 13.05

# Populate the newly inserted row 16 with the new "VINCENT MANUFACTURING" record.
$ws.Cells.Item(16, 1).Value = "VINCENT MANUFACTURING"
$ws.Cells.Item(16, 2).Value = "House Account"
$ws.Cells.Item(16, 3).Value = "030"
$ws.Cells.Item(16, 5).Value = "0008335"

# Column F is unused but every other row still carries an (empty, unstyled)
# cell placeholder for it; touch it (with a no-op default-font set) so it
# materializes the same way.
$ws.Cells.Item(16, 6).Font.Name = "Arial"
